# Modification des données des entreprise + création des données des stage
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# idEntreprise (column H) becomes a per-group sequential counter (1..19 / 1..11 / 1..19)
$idEntreprise = @{
    2 = 1;  3 = 2;  4 = 3;  5 = 4;  6 = 5;  7 = 6;  8 = 7;  9 = 8;  10 = 9;
    11 = 10; 12 = 11; 13 = 12; 14 = 13; 15 = 14; 16 = 15; 17 = 16; 18 = 17; 19 = 18; 20 = 19;
    21 = 1;  22 = 2;  23 = 3;  24 = 4;  25 = 5;  26 = 6;  27 = 7;  28 = 8;  29 = 9; 30 = 10; 31 = 11;
    32 = 1;  33 = 2;  34 = 3;  35 = 4;  36 = 5;  37 = 6;  38 = 7;  39 = 8;  40 = 9; 41 = 10; 42 = 11;
    43 = 12; 44 = 13; 45 = 14; 46 = 15; 47 = 16; 48 = 17; 49 = 18; 50 = 19
}

foreach ($row in $idEntreprise.Keys) {
    $ws.Cells.Item($row, 8).Value = $idEntreprise[$row]
}

# Update the absolute path recorded for the last opened folder
$wb.Path = "C:\Users\poifr1532452\Desktop\Projet\Sprint 1\Global\BD\"

# Update the saved view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H32").Select()
